$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for Wins / Losses / Ties (columns AD, AE, AF)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header formatting (bold, centered, bordered) from an existing
# header cell onto the new header cells
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill season record (Wins=77, Losses=85, Ties=0) for each data row (2-45)
for ($r = 2; $r -le 45; $r++) {
    $ws.Cells.Item($r, 30).Value = 77
    $ws.Cells.Item($r, 31).Value = 85
    $ws.Cells.Item($r, 32).Value = 0
}
